$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.309165477752686
$ws.Range("B1").Value = 3.289761781692505
$ws.Range("C1").Value = 2.600979328155518
$ws.Range("D1").Value = 1.358297944068909
$ws.Range("E1").Value = 1.00533127784729
